# Re-theme the deck's slide master from the "Integral" (Red Violet) theme to
# the stock "Office Theme" colour scheme.
#
# ppt/theme/theme1.xml (bound to the one SlideMaster in this deck) held the
# "Integral" / "Red Violet" clrScheme; the edit swaps its 12 theme colours
# for the default Office palette (the font scheme / format scheme were
# already identical between the two themes, so the colours are the only
# real content change reachable from the slide master).

function RgbValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p  = $ppt.ActivePresentation
$m  = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

# COM ColorScheme.Item ordering: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $cs.Item($i + 1).RGB = RgbValue $officeColors[$i]
}
